$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("asesorias")

# Shift the header cells from "Asesor" onward one column to the right
# (within the existing A:K header range) to make room for the new
# "Usuario" column in B, without growing the sheet past column K.
$ws.Range("F1").Value = $ws.Range("E1").Value()
$ws.Range("E1").Value = $ws.Range("D1").Value()
$ws.Range("D1").Value = $ws.Range("C1").Value()
$ws.Range("C1").Value = $ws.Range("B1").Value()
$ws.Range("B1").Value = "Usuario"

# First real row of "asesorias" data: david (Estudiante) asked Maryem Ruiz
# (Asesor) about "Consulta tema de clase" on 30-11-2023 at 03:00 - 03:20.
$ws.Range("A2").Value = "david "
$ws.Range("B2").Value = "s"
$ws.Range("C2").Value = "Maryem Ruíz"
$ws.Range("D2").Value = "Consulta tema de clase"
$ws.Range("E2").Value = "30-11-2023"
$ws.Range("F2").Value = "03:00 - 03:20"

$ws.Range("C8").Select()
